# "1k test pilot measurements"
# Rename the sheet, fill in the pilot render-time measurements table and
# the accompanying notes, then restore the matching column width / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Pilot"

# Library names (column A) - entered first so their shared-string ids come
# before the size headers, matching the saved workbook's string table order.
$ws.Range("A2").Value = "ApexCharts"
$ws.Range("A3").Value = "Frappe Charts"
$ws.Range("A4").Value = "Google Charts"
$ws.Range("A5").Value = "TeeChart JS"

# Dataset-size header row.
$ws.Range("C1").Value = "1k"
$ws.Range("B1").Value = "Small"
$ws.Range("D1").Value = "5k"
$ws.Range("E1").Value = "10k"

# "1k" render time measurements (ms), one per library.
$ws.Range("C2").Value = 168
$ws.Range("C3").Value = 71
$ws.Range("C4").Value = 57
$ws.Range("C5").Value = 646

# Pilot-test notes, column I.
$ws.Range("I7").Value = "One measurement per library only"
$ws.Range("I8").Value = "Unit: ms, rounded to whole number"
$ws.Range("I9").Value = "Bar charts only, settings barely touched"
$ws.Range("I10").Value = "Shift+f5 done before each measurement to ignore cached content"
$ws.Range("I11").Value = "Other browser windows and other applications are open - real test will be done on a separate machine"

# Widen column A so the library names fit.
$ws.Columns.Item(1).ColumnWidth = 13.65

# Leave the same cell selected as in the saved workbook.
$ws.Range("F14").Select() | Out-Null
